$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Títulos")

# Row 75 used to hold a distinct "2X VENCEDOR DA SUPERTAÇA DA ITÁLIA" title.
# It is unified with the text already used in row 74.
$ws.Range("A75").Value = "VENCEDOR DA SUPERTAÇA DA ITÁLIA"

# Rows 128-136 used to hold a distinct "10X CAMPEÃO DE ESPANHA" title.
# It is unified with the text already used in row 127.
$ws.Range("A128:A136").Value = "CAMPEÃO DE ESPANHA"

# Reflect the cell range that was selected while making the edit.
$ws.Range("A74:A75").Select()

$wb.Save()
